# vault backup: 2023-08-19 00:05:59
# Rebuild the "survey" table on Sheet1: a 5-column x 2-row literature-review
# table (Title / publication / Content / Advantage / shortcoming header,
# then one paper's data), with header/content formatting, column widths,
# row heights, a print setup, and the new cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Cell values (typed in this order so the shared-string table comes
#    out in the same order as the source document: header L-to-R, then
#    B2/A2, then D2/C2/E2).
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "publication"
$ws.Range("C1").Value = "Content"
$ws.Range("D1").Value = "Advantage"
$ws.Range("E1").Value = "shortcoming"

$ws.Range("B2").Value = "TMC"
$ws.Range("A2").Value = "Multipath Cooperative Routing with Efficient Acknowledgement for LEO Satellite Networks"

$ws.Range("D2").Value = "1. 减少ACK排队对数据传输的影响`n2.多径提高系统的吞度量"

$ws.Range("C2").Value = "1.多径路由+网络编码`n2. challenge:由于卫星链路的传播时延长，ACK的停等机制和hop-by-hop的ACK机制造成了长的时延`n3.多径选路：定时发送probe，然后在终点选取n条不相交的传输路径，然后传输给源节点`n4.网络编码：把多个数据包线性编码为一个batch`n5.NO-wait-ACK: 发送完当前batch后持续发送下一个batch,无需等待ACK，如果发送过的batch超时未收到ACK则重传`n"

$ws.Range("E2").Value = "1.选路策略的overhead过大`n2. 仅采用线性的网络编码"

# ---------------------------------------------------------------------
# 2. Formatting, applied in the order the styles first appear so the
#    generated style table lines up with the source file.
# ---------------------------------------------------------------------

# B2 ("TMC"): centered both ways, default font.
$ws.Range("B2").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B2").VerticalAlignment = -4108     # xlCenter

# A1 ("Title" header): big (18pt) font, centered both ways.
$ws.Range("A1").Font.Size = 18
$ws.Range("A1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1").VerticalAlignment = -4108     # xlCenter

# A2 (paper title): default font, vertically centered.
$ws.Range("A2").VerticalAlignment = -4108     # xlCenter

# B1 ("publication" header): big (18pt) font, vertically centered.
$ws.Range("B1").Font.Size = 18
$ws.Range("B1").Font.Family = 3
$ws.Range("B1").VerticalAlignment = -4108     # xlCenter

# C1:E1 ("Content"/"Advantage"/"shortcoming" headers): big (18pt) font,
# centered horizontally, top-aligned vertically.
$ws.Range("C1:E1").Font.Size = 18
$ws.Range("C1:E1").Font.Family = 3
$ws.Range("C1:E1").HorizontalAlignment = -4108 # xlCenter
$ws.Range("C1:E1").VerticalAlignment = -4160   # xlTop

# C2:E2 (content/advantage/shortcoming text): default font, left/top,
# wrapped text so the long notes render as a block.
$ws.Range("C2:E2").HorizontalAlignment = -4131 # xlLeft
$ws.Range("C2:E2").VerticalAlignment = -4160   # xlTop
$ws.Range("C2:E2").WrapText = $true

# ---------------------------------------------------------------------
# 3. Row heights / column widths.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 22.5
$ws.Rows.Item(2).RowHeight = 159.5

$ws.Columns.Item(1).ColumnWidth = 44.71
$ws.Columns.Item(2).ColumnWidth = 23.43
$ws.Columns.Item(3).ColumnWidth = 45
$ws.Columns.Item(4).ColumnWidth = 46.71
$ws.Columns.Item(5).ColumnWidth = 58.43

# ---------------------------------------------------------------------
# 4. Page setup / selection, matching the source's printer & view state.
# ---------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9    # xlPaperA4
$ws.PageSetup.Orientation = 1  # xlPortrait

$ws.Range("D6").Select()
